$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.501.82"
$ws.Range("E2").Value = "  +3.64%  "

$ws.Range("D3").Value = "1.591.61"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").Value = "  +0.93%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.47%  "

$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("E10").Value = "  +0.10%  "

$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("D12").Value = "1.818.51"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").Value = "1.610.40"
$ws.Range("E13").Value = "  +2.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.531"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "

$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").Value = "28.505.20"
$ws.Range("E16").Value = "  +3.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "0.0₃0708"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("E20").Value = "  -0.99%  "

$ws.Range("E21").Value = "  +0.88%  "

$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("E23").Value = "  -0.91%  "

$ws.Range("E24").Value = "  +2.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("E30").Value = "  -0.59%  "

$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("D34").Value = "1.403.31"
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("E36").Value = "  -9.77%  "

$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.30%  "

$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "

$ws.Range("E44").Value = "  +1.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.981"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.93%  "

$ws.Range("D47").Value = "1.728.54"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0522"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
